$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04110066666666667
$ws.Range("H2").Value = 0.123302
$ws.Range("I2").Value = 0.02671259512010182
$ws.Range("J2").Value = 0.02671259512010182
$ws.Range("M2").Value = 62.58874
$ws.Range("N2").Value = 187.76622
$ws.Range("O2").Value = 0.4331197020873656
$ws.Range("P2").Value = 0.4331197020873655
$ws.Range("Q2").Value = 2.572438939826667
$ws.Range("R2").Value = 23.15195045844
$ws.Range("S2").Value = 0.01156975124039892
$ws.Range("T2").Value = 0.01156975124039892

$ws.Range("G3").Value = 0.04110066666666667
$ws.Range("H3").Value = 0.123302
$ws.Range("I3").Value = 0.02671259512010182
$ws.Range("J3").Value = 0.02671259512010182
$ws.Range("O3").Value = 0.03956530671562308
$ws.Range("P3").Value = 0.03956530671562306
$ws.Range("Q3").Value = 0.2349912395371111
$ws.Range("R3").Value = 2.114921155834
$ws.Range("S3").Value = 0.001056892019097085
$ws.Range("T3").Value = 0.001056892019097084

$ws.Range("G4").Value = 0.04110066666666667
$ws.Range("H4").Value = 0.123302
$ws.Range("I4").Value = 0.02671259512010182
$ws.Range("J4").Value = 0.02671259512010182
$ws.Range("M4").Value = 16.124321
$ws.Range("N4").Value = 48.372963
$ws.Range("O4").Value = 0.111581749494894
$ws.Range("P4").Value = 0.111581749494894
$ws.Range("Q4").Value = 0.6627203426473333
$ws.Range("R4").Value = 5.964483083826
$ws.Range("S4").Value = 0.002980638097049729
$ws.Range("T4").Value = 0.002980638097049729

$ws.Range("G5").Value = 0.04110066666666667
$ws.Range("H5").Value = 0.123302
$ws.Range("I5").Value = 0.02671259512010182
$ws.Range("J5").Value = 0.02671259512010182
$ws.Range("M5").Value = 11.37633566666667
$ws.Range("N5").Value = 34.129007
$ws.Range("O5").Value = 0.07872526455705194
$ws.Range("P5").Value = 0.07872526455705191
$ws.Range("Q5").Value = 0.4675749801237779
$ws.Range("R5").Value = 4.208174821114
$ws.Range("S5").Value = 0.002102956117835431
$ws.Range("T5").Value = 0.00210295611783543

$ws.Range("G6").Value = 0.04110066666666667
$ws.Range("H6").Value = 0.123302
$ws.Range("I6").Value = 0.02671259512010182
$ws.Range("J6").Value = 0.02671259512010182
$ws.Range("M6").Value = 36.54706633333333
$ws.Range("N6").Value = 109.641199
$ws.Range("O6").Value = 0.2529089814311731
$ws.Range("P6").Value = 0.2529089814311731
$ws.Range("Q6").Value = 1.502108791010889
$ws.Range("R6").Value = 13.518979119098
$ws.Range("S6").Value = 0.006755855223208276
$ws.Range("T6").Value = 0.006755855223208276

$ws.Range("G7").Value = 0.04110066666666667
$ws.Range("H7").Value = 0.123302
$ws.Range("I7").Value = 0.02671259512010182
$ws.Range("J7").Value = 0.02671259512010182
$ws.Range("M7").Value = 12.15287633333333
$ws.Range("N7").Value = 36.458629
$ws.Range("O7").Value = 0.08409899571389245
$ws.Range("P7").Value = 0.08409899571389244
$ws.Range("Q7").Value = 0.4994913192175556
$ws.Range("R7").Value = 4.495421872958
$ws.Range("S7").Value = 0.002246502422512388
$ws.Range("T7").Value = 0.002246502422512388

$ws.Range("I8").Value = 0.6362555311831452
$ws.Range("J8").Value = 0.636255531183145
$ws.Range("M8").Value = 62.58874
$ws.Range("N8").Value = 187.76622
$ws.Range("O8").Value = 0.4331197020873656
$ws.Range("P8").Value = 0.4331197020873655
$ws.Range("Q8").Value = 61.27178945874667
$ws.Range("R8").Value = 551.44610512872
$ws.Range("S8").Value = 0.2755748061174824
$ws.Range("T8").Value = 0.2755748061174823

$ws.Range("I9").Value = 0.6362555311831452
$ws.Range("J9").Value = 0.636255531183145
$ws.Range("O9").Value = 0.03956530671562308
$ws.Range("P9").Value = 0.03956530671562306
$ws.Range("S9").Value = 0.02517364524077282
$ws.Range("T9").Value = 0.02517364524077281

$ws.Range("I10").Value = 0.6362555311831452
$ws.Range("J10").Value = 0.636255531183145
$ws.Range("M10").Value = 16.124321
$ws.Range("N10").Value = 48.372963
$ws.Range("O10").Value = 0.111581749494894
$ws.Range("P10").Value = 0.111581749494894
$ws.Range("Q10").Value = 15.78504378706533
$ws.Range("R10").Value = 142.065394083588
$ws.Range("S10").Value = 0.07099450529521842
$ws.Range("T10").Value = 0.07099450529521839

$ws.Range("I11").Value = 0.6362555311831452
$ws.Range("J11").Value = 0.636255531183145
$ws.Range("M11").Value = 11.37633566666667
$ws.Range("N11").Value = 34.129007
$ws.Range("O11").Value = 0.07872526455705194
$ws.Range("P11").Value = 0.07872526455705191
$ws.Range("Q11").Value = 11.13696239579244
$ws.Range("R11").Value = 100.232661562132
$ws.Range("S11").Value = 0.05008938501828071
$ws.Range("T11").Value = 0.05008938501828068

$ws.Range("I12").Value = 0.6362555311831452
$ws.Range("J12").Value = 0.636255531183145
$ws.Range("M12").Value = 36.54706633333333
$ws.Range("N12").Value = 109.641199
$ws.Range("O12").Value = 0.2529089814311731
$ws.Range("P12").Value = 0.2529089814311731
$ws.Range("Q12").Value = 35.77806732825822
$ws.Range("R12").Value = 322.002605954324
$ws.Range("S12").Value = 0.1609147383214792
$ws.Range("T12").Value = 0.1609147383214792

$ws.Range("I13").Value = 0.6362555311831452
$ws.Range("J13").Value = 0.636255531183145
$ws.Range("M13").Value = 12.15287633333333
$ws.Range("N13").Value = 36.458629
$ws.Range("O13").Value = 0.08409899571389245
$ws.Range("P13").Value = 0.08409899571389244
$ws.Range("Q13").Value = 11.89716361144489
$ws.Range("R13").Value = 107.074472503004
$ws.Range("S13").Value = 0.05350845118991169
$ws.Range("T13").Value = 0.05350845118991168

$ws.Range("G14").Value = 0.5185656666666667
$ws.Range("H14").Value = 1.555697
$ws.Range("I14").Value = 0.3370318736967531
$ws.Range("J14").Value = 0.3370318736967531
$ws.Range("M14").Value = 62.58874
$ws.Range("N14").Value = 187.76622
$ws.Range("O14").Value = 0.4331197020873656
$ws.Range("P14").Value = 0.4331197020873655
$ws.Range("Q14").Value = 32.45637168392667
$ws.Range("R14").Value = 292.10734515534
$ws.Range("S14").Value = 0.1459751447294843
$ws.Range("T14").Value = 0.1459751447294843

$ws.Range("G15").Value = 0.5185656666666667
$ws.Range("H15").Value = 1.555697
$ws.Range("I15").Value = 0.3370318736967531
$ws.Range("J15").Value = 0.3370318736967531
$ws.Range("O15").Value = 0.03956530671562308
$ws.Range("P15").Value = 0.03956530671562306
$ws.Range("Q15").Value = 2.964876209422111
$ws.Range("R15").Value = 26.683885884799
$ws.Range("S15").Value = 0.01333476945575317
$ws.Range("T15").Value = 0.01333476945575317

$ws.Range("G16").Value = 0.5185656666666667
$ws.Range("H16").Value = 1.555697
$ws.Range("I16").Value = 0.3370318736967531
$ws.Range("J16").Value = 0.3370318736967531
$ws.Range("M16").Value = 16.124321
$ws.Range("N16").Value = 48.372963
$ws.Range("O16").Value = 0.111581749494894
$ws.Range("P16").Value = 0.111581749494894
$ws.Range("Q16").Value = 8.361519268912334
$ws.Range("R16").Value = 75.253673420211
$ws.Range("S16").Value = 0.03760660610262585
$ws.Range("T16").Value = 0.03760660610262585

$ws.Range("G17").Value = 0.5185656666666667
$ws.Range("H17").Value = 1.555697
$ws.Range("I17").Value = 0.3370318736967531
$ws.Range("J17").Value = 0.3370318736967531
$ws.Range("M17").Value = 11.37633566666667
$ws.Range("N17").Value = 34.129007
$ws.Range("O17").Value = 0.07872526455705194
$ws.Range("P17").Value = 0.07872526455705191
$ws.Range("Q17").Value = 5.899377089208778
$ws.Range("R17").Value = 53.094393802879
$ws.Range("S17").Value = 0.0265329234209358
$ws.Range("T17").Value = 0.02653292342093579

$ws.Range("G18").Value = 0.5185656666666667
$ws.Range("H18").Value = 1.555697
$ws.Range("I18").Value = 0.3370318736967531
$ws.Range("J18").Value = 0.3370318736967531
$ws.Range("M18").Value = 36.54706633333333
$ws.Range("N18").Value = 109.641199
$ws.Range("O18").Value = 0.2529089814311731
$ws.Range("P18").Value = 0.2529089814311731
$ws.Range("Q18").Value = 18.95205381785589
$ws.Range("R18").Value = 170.568484360703
$ws.Range("S18").Value = 0.08523838788648559
$ws.Range("T18").Value = 0.08523838788648559

$ws.Range("G19").Value = 0.5185656666666667
$ws.Range("H19").Value = 1.555697
$ws.Range("I19").Value = 0.3370318736967531
$ws.Range("J19").Value = 0.3370318736967531
$ws.Range("M19").Value = 12.15287633333333
$ws.Range("N19").Value = 36.458629
$ws.Range("O19").Value = 0.08409899571389245
$ws.Range("P19").Value = 0.08409899571389244
$ws.Range("Q19").Value = 6.302064417712556
$ws.Range("R19").Value = 56.718579759413
$ws.Range("S19").Value = 0.02834404210146838
$ws.Range("T19").Value = 0.02834404210146838
